$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: D = Price column (text-typed, numeric-looking values must
# stay text so they keep their original inlineStr semantics), E = Volume(1h).
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '27.662.67'
$ws.Range("E2").Value = '  +1.06%  '
Set-TextValue $ws.Range("D3") '1.872.42'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  +0.23%  '
Set-TextValue $ws.Range("D5") '331.65'
$ws.Range("E5").Value = '  +2.49%  '
Set-TextValue $ws.Range("D7") '0.4716'
$ws.Range("E7").Value = '  +4.15%  '
Set-TextValue $ws.Range("D8") '0.3940'
$ws.Range("E8").Value = '  +1.87%  '
Set-TextValue $ws.Range("D9") '47.93'
$ws.Range("E9").Value = '  -1.91%  '
Set-TextValue $ws.Range("D10") '0.08024'
$ws.Range("E10").Value = '  +1.52%  '
Set-TextValue $ws.Range("D11") '1.026'
$ws.Range("E11").Value = '  +0.89%  '
Set-TextValue $ws.Range("D12") '22.02'
$ws.Range("E12").Value = '  +3.06%  '
Set-TextValue $ws.Range("D13") '1.851.58'
$ws.Range("E13").Value = '  -0.18%  '
Set-TextValue $ws.Range("D14") '5.956'
$ws.Range("E14").Value = '  +0.53%  '
Set-TextValue $ws.Range("D15") '7.116'
Set-TextValue $ws.Range("D16") '1.006'
$ws.Range("E16").Value = '  +0.45%  '
Set-TextValue $ws.Range("D17") '0.00001047'
$ws.Range("E17").Value = '  +1.48%  '
Set-TextValue $ws.Range("D18") '86.91'
$ws.Range("E18").Value = '  +1.19%  '
Set-TextValue $ws.Range("D19") '0.06668'
$ws.Range("E19").Value = '  +2.30%  '
Set-TextValue $ws.Range("D20") '17.21'
$ws.Range("E20").Value = '  +1.00%  '
Set-TextValue $ws.Range("D22") '27.667.34'
$ws.Range("E22").Value = '  +1.08%  '
Set-TextValue $ws.Range("D23") '5.510'
$ws.Range("E23").Value = '  -0.31%  '
Set-TextValue $ws.Range("D24") '10.95'
$ws.Range("E24").Value = '  +1.17%  '
Set-TextValue $ws.Range("D25") '2.305'
$ws.Range("E25").Value = '  +0.91%  '
Set-TextValue $ws.Range("D26") '2.076.52'
$ws.Range("E26").Value = '  +0.04%  '
Set-TextValue $ws.Range("D27") '158.47'
$ws.Range("E27").Value = '  +3.02%  '
Set-TextValue $ws.Range("D28") '20.14'
$ws.Range("E28").Value = '  +1.83%  '
Set-TextValue $ws.Range("D29") '2.100'
$ws.Range("E29").Value = '  +1.22%  '
Set-TextValue $ws.Range("D30") '5.563'
$ws.Range("E30").Value = '  +2.27%  '
Set-TextValue $ws.Range("D31") '122.15'
$ws.Range("E31").Value = '  +0.95%  '
Set-TextValue $ws.Range("D32") '0.9729'
$ws.Range("E32").Value = '  +3.83%  '
Set-TextValue $ws.Range("D33") '0.09524'
$ws.Range("E33").Value = '  +2.57%  '
Set-TextValue $ws.Range("D34") '1.445'
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("E35").Value = '  -0.22%  '
Set-TextValue $ws.Range("D36") '5.328'
$ws.Range("E36").Value = '  +1.47%  '
Set-TextValue $ws.Range("D37") '0.06095'
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("E38").Value = '  +0.65%  '
Set-TextValue $ws.Range("D39") '1.228'
$ws.Range("E39").Value = '  +0.61%  '
Set-TextValue $ws.Range("D40") '8.151'
$ws.Range("E40").Value = '  -0.69%  '
Set-TextValue $ws.Range("D41") '0.6014'
$ws.Range("E41").Value = '  +1.81%  '
Set-TextValue $ws.Range("D42") '0.1904'
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("E43").Value = '  +1.37%  '
Set-TextValue $ws.Range("D44") '1.252'
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("E45").Value = '  +1.46%  '
Set-TextValue $ws.Range("D46") '12.16'
$ws.Range("E46").Value = '  +1.63%  '
Set-TextValue $ws.Range("D47") '1.941'
$ws.Range("E47").Value = '  +0.96%  '
Set-TextValue $ws.Range("D48") '3.382'
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("E49").Value = '  +1.77%  '
Set-TextValue $ws.Range("D50") '115.05'
$ws.Range("E50").Value = '  +6.31%  '
$ws.Range("E51").Value = '  +8.39%  '
